$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header row gets new "설명" columns interleaved, plus two more photo slots ---
$ws.Range("F1").Value = "사진1설명"
$ws.Range("G1").Value = "사진2"
$ws.Range("H1").Value = "사진2설명"
$ws.Range("I1").Value = "사진3"
$ws.Range("J1").Value = "사진3설명"
$ws.Range("K1").Value = "사진4"
$ws.Range("L1").Value = "사진4설명"
$ws.Range("M1").Value = "사진5"
$ws.Range("N1").Value = "사진5설명"

# --- Rows 2-4: new empty-string (text) cells for columns H through N ---
# A direct assignment of "" is treated (as in real Excel COM) as clearing
# the cell rather than storing an empty string, so use a formula that
# evaluates to an empty string to materialize a Text-typed blank cell.
$ws.Range("H2:N4").Formula = '=""'

# --- Text corrections in rows 3 and 4 ---
$ws.Range("C3").Value = "브레이크 레버 간격 과다"
$ws.Range("D3").Value = "레버 간격 조정 후 점검"
$ws.Range("C4").Value = "시동 경고등 점등"
